$d = $word.ActiveDocument

# Locate the document title - "The World's Best Selling Cars" - which is the
# Heading1-styled paragraph at the top of the document, then apply yellow
# highlighting to it (matching <w:highlight w:val="yellow"/> added to every
# run in that paragraph).
$titleRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Style.NameLocal -eq "Heading 1") {
        $titleRange = $para.Range
        break
    }
}

if ($titleRange -eq $null) {
    # Fallback: search for the title text directly.
    $titleRange = $d.Content
    $titleRange.Find.Execute("The World")
    $titleRange = $titleRange.Paragraphs(1).Range
}

$titleRange.HighlightColorIndex = 7

Write-Host "Highlighted title:" $titleRange.Text
